$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "29.578.91"
$ws.Cells.Item(2, 5).Value = "  +2.49%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.991.46"
$ws.Cells.Item(3, 5).Value = "  +6.04%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.17%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "329.71"
$ws.Cells.Item(5, 5).Value = "  +1.60%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.11%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4681"
$ws.Cells.Item(7, 5).Value = "  +1.44%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3955"
$ws.Cells.Item(8, 5).Value = "  +2.17%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "46.52"
$ws.Cells.Item(9, 5).Value = "  -0.29%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.08102"
$ws.Cells.Item(10, 5).Value = "  +3.28%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.001"
$ws.Cells.Item(11, 5).Value = "  +1.81%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "22.92"
$ws.Cells.Item(12, 5).Value = "  +5.51%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.992.83"
$ws.Cells.Item(13, 5).Value = "  +5.35%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.248"
$ws.Cells.Item(14, 5).Value = "  +3.65%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.871"
$ws.Cells.Item(15, 5).Value = "  +3.77%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.07127"
$ws.Cells.Item(16, 5).Value = "  +2.15%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "88.83"
$ws.Cells.Item(17, 5).Value = "  +0.74%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "1.003"
$ws.Cells.Item(18, 5).Value = "  -0.09%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.00001004"
$ws.Cells.Item(19, 5).Value = "  +0.91%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "17.44"
$ws.Cells.Item(20, 5).Value = "  +3.27%  "
$ws.Cells.Item(21, 5).Value = "  -0.12%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "29.601.48"
$ws.Cells.Item(22, 5).Value = "  +2.60%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.559"
$ws.Cells.Item(23, 5).Value = "  +5.79%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.26"
$ws.Cells.Item(24, 5).Value = "  +2.84%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.115"
$ws.Cells.Item(25, 5).Value = "  +0.54%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "157.93"
$ws.Cells.Item(27, 5).Value = "  +2.11%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "5.990"
$ws.Cells.Item(28, 5).Value = "  +0.44%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "120.21"
$ws.Cells.Item(29, 5).Value = "  +2.21%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.946"
$ws.Cells.Item(30, 5).Value = "  +2.06%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.09459"
$ws.Cells.Item(31, 5).Value = "  +1.17%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.9162"
$ws.Cells.Item(32, 5).Value = "  +1.83%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.281"
$ws.Cells.Item(33, 5).Value = "  +0.35%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.355"
$ws.Cells.Item(34, 5).Value = "  +2.86%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.186"
$ws.Cells.Item(35, 5).Value = "  -2.05%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.05857"
$ws.Cells.Item(36, 5).Value = "  +2.25%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.176"
$ws.Cells.Item(37, 5).Value = "  +0.70%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02131"
$ws.Cells.Item(38, 5).Value = "  +2.86%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.000003338"
$ws.Cells.Item(39, 5).Value = "  +86.68%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "7.912"
$ws.Cells.Item(40, 5).Value = "  +3.71%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.5790"
$ws.Cells.Item(41, 5).Value = "  +2.51%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1828"
$ws.Cells.Item(42, 5).Value = "  +3.37%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "9.908"
$ws.Cells.Item(43, 5).Value = "  +2.17%  "
$ws.Cells.Item(44, 2).Value = "MXToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.799"
$ws.Cells.Item(44, 5).Value = "  +10.47%  "
$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "12.04"
$ws.Cells.Item(45, 5).Value = "  +0.75%  "
$ws.Cells.Item(46, 5).Value = "  +1.29%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.210"
$ws.Cells.Item(47, 5).Value = "  -1.10%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.873"
$ws.Cells.Item(48, 5).Value = "  +1.78%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.06968"
$ws.Cells.Item(49, 5).Value = "  -0.99%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "114.10"
$ws.Cells.Item(50, 5).Value = "  +1.44%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.3087"
$ws.Cells.Item(51, 5).Value = "  +7.98%  "
